# Auto-generated edit script: update crypto price/volume data per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "42.571.69"
$ws.Range("E2").Value = "  +0.28%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.282.57"
$ws.Range("E3").Value = "  -0.79%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.15%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.36"
$ws.Range("E5").Value = "  -1.30%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.08"
$ws.Range("E6").Value = "  +0.78%  "

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  -0.82%  "

# Row 8: USDC
$ws.Range("E8").Value = "  -0.16%  "

# Row 9: Cardano
$ws.Range("E9").Value = "  -1.30%  "

# Row 10: Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.57"
$ws.Range("E10").Value = "  -0.92%  "

# Row 11: Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0901"
$ws.Range("E11").Value = "  -0.69%  "

# Row 12: Polkadot
$ws.Range("E12").Value = "  +0.53%  "

# Row 13: TRON
$ws.Range("E13").Value = "  +2.51%  "

# Row 14: Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.996"
$ws.Range("E14").Value = "  +3.31%  "

# Row 15: Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.20"
$ws.Range("E15").Value = "  -0.86%  "

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "2.629.57"
$ws.Range("E16").Value = "  -0.77%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "2.270.59"
$ws.Range("E17").Value = "  -1.62%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "42.731.24"
$ws.Range("E18").Value = "  +0.70%  "

# Row 19: Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.37"
$ws.Range("E19").Value = "  -1.33%  "

# Row 20: ShibaInu
$ws.Range("E20").Value = "  -0.57%  "

# Row 21: InternetComputer(DFINITY)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.54"
$ws.Range("E21").Value = "  +21.50%  "

# Row 22: Litecoin
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.85"
$ws.Range("E22").Value = "  +0.70%  "

# Row 23: PancakeSwap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.53"

# Row 24: BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.16"
$ws.Range("E24").Value = "  -4.50%  "

# Row 25: ImmutableX
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("E25").Value = "  -2.74%  "

# Row 26: Dai
$ws.Range("E26").Value = "  +0.36%  "

# Row 27: Cosmos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.81"
$ws.Range("E27").Value = "  -0.13%  "

# Row 28: Filecoin
$ws.Range("E28").Value = "  +20.94%  "

# Row 29: Toncoin
$ws.Range("E29").Value = "  -0.24%  "

# Row 30: EthereumClassic
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.44"
$ws.Range("E30").Value = "  -1.42%  "

# Row 31: InjectiveProtocol
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.03"
$ws.Range("E31").Value = "  +3.99%  "

# Row 32: Monero
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.27"
$ws.Range("E32").Value = "  +0.60%  "

# Row 33: Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0871"

# Row 34: Stellar
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.131"
$ws.Range("E34").Value = "  -3.02%  "

# Row 35: WEMIXToken
$ws.Range("E35").Value = "  +1.22%  "

# Row 36: Kaspa
$ws.Range("E36").Value = "  -3.11%  "

# Row 37: RenderToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.53"
$ws.Range("E37").Value = "  -1.13%  "

# Row 38: VeChain
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0350"
$ws.Range("E38").Value = "  -4.71%  "

# Row 39: NEARProtocol
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.79"
$ws.Range("E39").Value = "  +1.61%  "

# Row 40: LidoDAOToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.65"
$ws.Range("E40").Value = "  -3.79%  "

# Row 41: ARBITRUM
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.57"
$ws.Range("E41").Value = "  +5.39%  "

# Row 42: Algorand -> MultiversX
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.30"
$ws.Range("E42").Value = "  +0.99%  "

# Row 43: MultiversX -> Algorand
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.231"
$ws.Range("E43").Value = "  +1.89%  "

# Row 44: FirstDigitalUSD -> BitcoinSV
$ws.Range("B44").Value = "BitcoinSV"
$ws.Range("C44").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "95.11"
$ws.Range("E44").Value = "  +0.92%  "

# Row 45: BitcoinSV -> FirstDigitalUSD
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.02%  "

# Row 46: Celestia
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.17"
$ws.Range("E46").Value = "  +0.98%  "

# Row 47: Maker -> Aave
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "113.62"
$ws.Range("E47").Value = "  +0.68%  "

# Row 48: Aave -> Maker
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.729.80"
$ws.Range("E48").Value = "  +8.81%  "

# Row 49: ordi
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "79.17"
$ws.Range("E49").Value = "  -2.98%  "

# Row 50: FraxShare
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.74"
$ws.Range("E50").Value = "  -1.94%  "

# Row 51: THORChain
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.18"
$ws.Range("E51").Value = "  +0.01%  "

